$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix truncated image URL: ".pn" -> ".png"
$ws.Range("A2").Value = "https://res.cloudinary.com/dwbpf1nax/image/upload/v1693738214/588-1_boo31q.png"

# Update date value in D2
$ws.Range("D2").Value = 45192

# Move selection from B8 to C9
$null = $ws.Range("C9").Select()
